$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

function Replace-TextInRange($range, $old, $new) {
    $ok = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace (scoped) failed for: $old"
    }
}

# Title
Replace-Text "Unraveling the Enigma of Consciousness: A Neuroscientific Expedition" "The Marvelous Interplay: An Exploration of Chemistry and Everyday Phenomena"

# Author name (collapses three runs "Dr" + "." + " Anya Chandra" into one run)
Replace-Text "Dr. Anya Chandra" "Alex Sterling"

# Email parts - the "com" replacement must stay scoped to the email paragraph only,
# otherwise it would also hit "com" inside "communicate"/"complexity" later in the doc.
Replace-Text "anyachandra@neurosciences" "alex_sterling@xyz"
$emailPara = $d.Paragraphs.Item(3).Range
Replace-TextInRange $emailPara "com" "academy"

# Intro paragraph sentences
Replace-Text "Consciousness, the enigmatic faculty that allows us to experience subjective awareness, has captivated thinkers since the dawn of civilization" "In the vast world of science, chemistry stands as a captivating subject that delves into the intricate interactions of matter at the microscopic level"
Replace-Text "As we navigate the 21st century, armed with advanced brain imaging and experimental tools, we stand poised to demystify this profound mystery" "Its principles are all around us, shaping our lives in countless ways"
Replace-Text "In this scientific odyssey, we will traverse the intricate neural landscapes of consciousness, exploring the biological underpinnings, cognitive processes, and evolutionary foundations that shape our sentient existence" "From the foods we eat and the clothes we wear to the medicines that heal our bodies and the products that power our daily routines, chemistry plays a vital role in the fabric of our existence"

# Second block
Replace-Text "Venturing into the labyrinthine depths of the human brain, contemporary neuroscience has begun to elucidate the neural correlates and processes underlying conscious experience" "Journey into the world of chemistry and discover the fascinating explanations behind the world's physical phenomena"
Replace-Text "From the enigmatic interplay of brain regions, including the cerebral cortex, thalamus, and brainstem, to the dynamic interactions of neuronal networks, scientists are teasing apart the mechanisms that orchestrate the symphony of consciousness" "Delve into the realm of elements and discover the role they play in creating the substances that make up the world around you"

# Insert the new sentence (". Witness the awe-inspiring ...") right after the sentence above, before its trailing period run.
$findRng = $d.Content
$ok = $findRng.Find.Execute("Delve into the realm of elements and discover the role they play in creating the substances that make up the world around you", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Find failed locating Delve-into sentence" }
$findRng.Collapse(0)
$findRng.InsertAfter(". Witness the awe-inspiring spectacle of chemical reactions as molecules dance and rearrange themselves, leading to the formation of entirely new substances and revealing the intricate tapestry that weaves our world together")

# Third block
Replace-Text "Furthermore, the exploration of consciousness through the lens of evolutionary biology offers a profound perspective on its origins and adaptive functions" "Uncover the fundamental principles that govern the behavior of matter"
Replace-Text "Consciousness likely arose through natural selection, granting our ancestors the ability to navigate an intricate environment, make informed decisions, and communicate with ever-increasing complexity" "Investigate the interactions between atoms and molecules, and learn how the laws of thermodynamics shape their transformations"
Replace-Text "This evolutionary journey, spanning eons, has left an imprint on the neural architecture of consciousness, revealing clues that can help unravel its intricate web" "Explore the concepts of acids and bases and delve into the fascinating world of chemical equilibrium, where reactions seek to achieve a delicate balance"

# Insert new sentence after "...achieve a delicate balance"
$findRng2 = $d.Content
$ok = $findRng2.Find.Execute("Explore the concepts of acids and bases and delve into the fascinating world of chemical equilibrium, where reactions seek to achieve a delicate balance", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Find failed locating Explore-concepts sentence" }
$findRng2.Collapse(0)
$findRng2.InsertAfter(". Discover the diverse forms of chemical energy and witness the majestic displays of energy transfer that power everything from the burning of fuels to the function of our own bodies")

# Summary paragraph
Replace-Text "Our quest to decipher the enigmatic puzzle of consciousness has led us through the intricate pathways of the human brain, from the neural underpinnings to the evolutionary foundations of subjective experience" "We embarked on a journey into the realm of chemistry, unraveling the mysteries of matter and its interactions"
Replace-Text "As we continue to navigate the uncharted waters of consciousness, armed with advancements in neuroscience and interdisciplinary perspectives, we can anticipate a growing understanding of this profound phenomenon--a testament to the indomitable spirit of exploration inherent in the human species" "From the fundamental principles that govern molecular behavior to the practical applications that touch every aspect of our lives, chemistry stands as a cornerstone of scientific understanding. By exploring the wonders of chemical reactions, we gain insights into the natural world and deepen our appreciation for the intricate processes that sustain life and shape our surroundings"

# Add a trailing empty paragraph at the end of the document
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
